$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1997.7778
$ws.Range("I40").Value = 1854.2858
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1854.2858
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1679.2858
$ws.Range("N40").Value = -2850

$ws.Range("H55").Value = 541.1875
$ws.Range("I55").Value = 132.2
$ws.Range("J55").Value = 727.0909
$ws.Range("K55").Value = 132.2
$ws.Range("L55").Value = 727.0909
$ws.Range("M55").Value = 81.80000000000001
$ws.Range("N55").Value = -1155.0909

$ws.Range("H97").Value = 6973.3335
$ws.Range("J97").Value = 6250
$ws.Range("L97").Value = 18750
$ws.Range("N97").Value = -19742

$ws.Range("H134").Value = 73333.336
$ws.Range("J134").Value = 73333.336
$ws.Range("L134").Value = 73333.336
$ws.Range("N134").Value = -83473.336

$ws.Range("H137").Value = 2935.8845
$ws.Range("I137").Value = 1314.1
$ws.Range("K137").Value = 3942.3
$ws.Range("M137").Value = -1392.3

$ws.Range("H138").Value = 5070.9287
$ws.Range("I138").Value = 4698.25
$ws.Range("J138").Value = 5220
$ws.Range("K138").Value = 14094.75
$ws.Range("L138").Value = 15660
$ws.Range("M138").Value = -8954.75
$ws.Range("N138").Value = -25940

$ws.Range("H139").Value = 100000
$ws.Range("I139").Value = 100000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 100000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -94860
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1273
$ws.Range("I2").Value = 1162.2
$ws.Range("K2").Value = 1162.2
$ws.Range("M2").Value = -1049.2

$ws.Range("H40").Value = 17206
$ws.Range("J40").Value = 17206
$ws.Range("L40").Value = 17206
$ws.Range("N40").Value = -17558

$ws.Range("H116").Value = 1273
$ws.Range("I116").Value = 1162.2
$ws.Range("K116").Value = 1162.2
$ws.Range("M116").Value = 1131.8

$ws.Range("H130").Value = 24166.334
$ws.Range("J130").Value = 24166.334
$ws.Range("L130").Value = 24166.334
$ws.Range("N130").Value = -34206.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1273
$ws.Range("I3").Value = 1162.2
$ws.Range("K3").Value = 1162.2
$ws.Range("M3").Value = -1048.2

$ws.Range("H64").Value = 1044.6923
$ws.Range("I64").Value = 871
$ws.Range("J64").Value = 1193.5714
$ws.Range("K64").Value = 871
$ws.Range("L64").Value = 1193.5714
$ws.Range("M64").Value = -646
$ws.Range("N64").Value = -1643.5714

$ws.Range("H67").Value = 1044.6923
$ws.Range("I67").Value = 871
$ws.Range("J67").Value = 1193.5714
$ws.Range("K67").Value = 871
$ws.Range("L67").Value = 1193.5714
$ws.Range("M67").Value = -91
$ws.Range("N67").Value = -2753.5714

$ws.Range("H99").Value = 1328.125
$ws.Range("I99").Value = 1328.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1328.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 169.875
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11108.4
$ws.Range("I99").Value = 8299.286
$ws.Range("K99").Value = 8299.286
$ws.Range("M99").Value = -6801.286

$ws.Range("H126").Value = 11108.4
$ws.Range("I126").Value = 8299.286
$ws.Range("K126").Value = 24897.858
$ws.Range("M126").Value = -22427.858

$ws.Range("H134").Value = 1783.3422
$ws.Range("I134").Value = 1219.2258
$ws.Range("J134").Value = 4281.5713
$ws.Range("K134").Value = 3657.6774
$ws.Range("L134").Value = 12844.7139
$ws.Range("M134").Value = -1122.6774
$ws.Range("N134").Value = -17914.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 30.928572
$ws.Range("J12").Value = 42.714287
$ws.Range("L12").Value = 128.142861
$ws.Range("N12").Value = -474.142861

$ws.Range("H14").Value = 1739.75
$ws.Range("I14").Value = 1739.75
$ws.Range("K14").Value = 5219.25
$ws.Range("M14").Value = -5046.25

$ws.Range("H22").Value = 5184.75
$ws.Range("J22").Value = 5184.75
$ws.Range("L22").Value = 15554.25
$ws.Range("N22").Value = -15892.25

$ws.Range("H27").Value = 5184.75
$ws.Range("J27").Value = 5184.75
$ws.Range("L27").Value = 15554.25
$ws.Range("N27").Value = -15758.25

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9796.4
$ws.Range("I80").Value = 9661
$ws.Range("J80").Value = 9999.5
$ws.Range("K80").Value = 9661
$ws.Range("L80").Value = 9999.5
$ws.Range("M80").Value = -8663
$ws.Range("N80").Value = -11995.5

$ws.Range("H83").Value = 9796.4
$ws.Range("I83").Value = 9661
$ws.Range("J83").Value = 9999.5
$ws.Range("K83").Value = 48305
$ws.Range("L83").Value = 49997.5
$ws.Range("M83").Value = -43313
$ws.Range("N83").Value = -59981.5

$ws.Range("H126").Value = 3343.7144
$ws.Range("I126").Value = 2482.4
$ws.Range("J126").Value = 3822.2222
$ws.Range("K126").Value = 7447.200000000001
$ws.Range("L126").Value = 11466.6666
$ws.Range("M126").Value = -4977.200000000001
$ws.Range("N126").Value = -16406.6666

$ws.Range("H132").Value = 2355.3103
$ws.Range("I132").Value = 1923
$ws.Range("K132").Value = 5769
$ws.Range("M132").Value = -3239

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 800.75
$ws.Range("J22").Value = 601
$ws.Range("L22").Value = 601
$ws.Range("N22").Value = -1191

$ws.Range("H27").Value = 800.75
$ws.Range("J27").Value = 601
$ws.Range("L27").Value = 601
$ws.Range("N27").Value = -815

$ws.Range("H46").Value = 3232.353
$ws.Range("I46").Value = 1995
$ws.Range("K46").Value = 1995
$ws.Range("M46").Value = -1807

$ws.Range("H93").Value = 1362.25
$ws.Range("I93").Value = 1332.75
$ws.Range("J93").Value = 1391.75
$ws.Range("K93").Value = 1332.75
$ws.Range("L93").Value = 1391.75
$ws.Range("M93").Value = -84.75
$ws.Range("N93").Value = -3887.75

$ws.Range("H132").Value = 2887.9697
$ws.Range("I132").Value = 1958.16
$ws.Range("J132").Value = 5793.625
$ws.Range("K132").Value = 5874.48
$ws.Range("L132").Value = 17380.875
$ws.Range("M132").Value = -3344.48
$ws.Range("N132").Value = -22440.875

$ws.Range("H136").Value = 2873.75
$ws.Range("I136").Value = 2873.75
$ws.Range("K136").Value = 8621.25
$ws.Range("M136").Value = -6071.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1255.8077
$ws.Range("I122").Value = 1206.04
$ws.Range("K122").Value = 3618.12
$ws.Range("M122").Value = -1168.12

$ws.Range("H126").Value = 93717.82000000001
$ws.Range("I126").Value = 126424.5
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 379273.5
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -376803.5
$ws.Range("N126").Value = -24440
